$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Low-grade glioma")
$ws.Range("C3").Value = 0.603839616038396
$ws.Range("C6").Value = 0.0754924507549245
$ws.Range("C7").Value = 0.218878112188781
$ws.Range("C8").Value = 0.789421057894211

$ws = $wb.Worksheets.Item("Ependymoma")
$ws.Range("C3").Value = 0.388761123887611
$ws.Range("C6").Value = 0.256674332566743
$ws.Range("C7").Value = 0.285171482851715
$ws.Range("C8").Value = 0.967303269673033

$ws = $wb.Worksheets.Item("DIPG or DMG")
$ws.Range("C3").Value = 0.615538446155384
$ws.Range("C6").Value = 0.0173982601739826
$ws.Range("C7").Value = 0.768523147685231

$ws = $wb.Worksheets.Item("ATRT")
$ws.Range("C6").Value = 0.4997500249975
$ws.Range("C7").Value = 0.275772422757724
$ws.Range("C8").Value = 0.173682631736826

$ws = $wb.Worksheets.Item("Other high-grade glioma")
$ws.Range("C3").Value = 0.565043495650435
$ws.Range("C6").Value = 0.8001199880012
$ws.Range("C8").Value = 0.655634436556344

$ws = $wb.Worksheets.Item("Meningioma")
$ws.Range("C3").Value = 0.946905309469053
$ws.Range("C6").Value = 0.689431056894311
$ws.Range("C7").Value = 0.318368163183682

$ws = $wb.Worksheets.Item("Neurofibroma plexiform")
$ws.Range("C4").Value = 0.0003999600039996
$ws.Range("C5").Value = 0.154284571542846
$ws.Range("C6").Value = 0.0600939906009399
$ws.Range("C7").Value = 0.292170782921708

$ws = $wb.Worksheets.Item("Oligodendroglioma")
$ws.Range("C3").Value = 0.195880411958804
$ws.Range("C4").Value = 0.0007999200079992
$ws.Range("C5").Value = 0.478752124787521
$ws.Range("C7").Value = 0.388561143885611

$ws = $wb.Worksheets.Item("Non-neoplastic tumor")
$ws.Range("C3").Value = 0.916808319168083
$ws.Range("C5").Value = 0.0001999800019998
$ws.Range("C6").Value = 0.255874412558744
$ws.Range("C7").Value = 0.855414458554145

$ws = $wb.Worksheets.Item("Mixed neuronal-glial tumor")
$ws.Range("C3").Value = 0.194980501949805
$ws.Range("C6").Value = 0.526747325267473
$ws.Range("C7").Value = 0.788321167883212
$ws.Range("C8").Value = 0.695630436956304

$ws = $wb.Worksheets.Item("Medulloblastoma")
$ws.Range("C3").Value = 0.783721627837216
$ws.Range("C6").Value = 0.465753424657534
$ws.Range("C7").Value = 0.295770422957704
$ws.Range("C8").Value = 0.674632536746325

$ws = $wb.Worksheets.Item("Schwannoma")
$ws.Range("C3").Value = 0.805919408059194
$ws.Range("C5").Value = 0.503049695030497
$ws.Range("C6").Value = 0.935206479352065
$ws.Range("C7").Value = 0.846415358464154

$ws = $wb.Worksheets.Item("Mesenchymal tumor")
$ws.Range("C3").Value = 0.121687831216878
$ws.Range("C5").Value = 0.0001999800019998
$ws.Range("C6").Value = 0.0223977602239776
$ws.Range("C7").Value = 0.849315068493151
$ws.Range("C8").Value = 0.488151184881512

$ws = $wb.Worksheets.Item("Germ cell tumor")
$ws.Range("C3").Value = 0.610838916108389
$ws.Range("C4").Value = 0.0002999700029997
$ws.Range("C5").Value = 0.0077992200779922
$ws.Range("C6").Value = 0.853914608539146
$ws.Range("C7").Value = 0.352564743525647
$ws.Range("C8").Value = 0.283271672832717

$ws = $wb.Worksheets.Item("Craniopharyngioma")
$ws.Range("C3").Value = 0.733826617338266
$ws.Range("C5").Value = 0.004999500049995
$ws.Range("C6").Value = 0.163983601639836
$ws.Range("C7").Value = 0.581041895810419

$ws = $wb.Worksheets.Item("Other tumor")
$ws.Range("C3").Value = 0.0472952704729527
$ws.Range("C5").Value = 0.0006999300069993
$ws.Range("C6").Value = 0.879212078792121
$ws.Range("C7").Value = 0.994100589941006

Write-Host "Updated p-values for derived cell lines -> solid tumors comparison across all histology sheets."